$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.347.19"
$ws.Range("E2").Value = "  -4.42%  "
$ws.Range("D3").Value = "1.762.43"
$ws.Range("E3").Value = "  -4.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  -2.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4272"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3620"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07047"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8310"
$ws.Range("E10").Value = "  -3.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.16"
$ws.Range("E11").Value = "  -2.48%  "
$ws.Range("D12").Value = "1.744.26"
$ws.Range("E12").Value = "  -5.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.236"
$ws.Range("E13").Value = "  -3.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.381"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06791"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.11"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008650"
$ws.Range("E18").Value = "  -2.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  -3.37%  "
$ws.Range("D21").Value = "26.218.75"
$ws.Range("E21").Value = "  -4.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.995"
$ws.Range("E22").Value = "  -3.25%  "
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").Value = "1.970.57"
$ws.Range("E24").Value = "  -4.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.906"
$ws.Range("E25").Value = "  -4.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.86"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.07"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.032"
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.43"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -8.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08874"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7236"
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.111"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.305"
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.715"
$ws.Range("E36").Value = "  -9.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.063"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05091"
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01879"
$ws.Range("E39").Value = "  -3.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1603"
$ws.Range("E40").Value = "  -3.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4889"
$ws.Range("E41").Value = "  -3.83%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.165"
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.476"
$ws.Range("E43").Value = "  -11.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.988"
$ws.Range("E44").Value = "  -4.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.57"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.04"
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06185"
$ws.Range("E48").Value = "  -4.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4457"
$ws.Range("E49").Value = "  -4.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.568"
$ws.Range("E50").Value = "  -3.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.716"
$ws.Range("E51").Value = "  -1.15%  "

Write-Host "Applied cryptos update"
